$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Adjust the time-of-day fraction on several "datetime" cells in column A ---
# These rows shift from 10:00 to 11:00 (and one from 11:00 back to 10:00),
# same calendar day, only the serial's fractional part changes.
$ws.Range("A4").Value  = 35156.45833333334
$ws.Range("A16").Value = 35521.45833333334
$ws.Range("A28").Value = 35886.45833333334
$ws.Range("A40").Value = 36251.45833333334
$ws.Range("A52").Value = 36617.45833333334
$ws.Range("A57").Value = 36770.41666666666
$ws.Range("A76").Value = 37347.45833333334
$ws.Range("A88").Value = 37712.45833333334
$ws.Range("A100").Value = 38078.45833333334
$ws.Range("A112").Value = 38443.45833333334
$ws.Range("A130").Value = 38991.45833333334

# --- Updated OHLC values (open/high/low/close all equal) for several rows ---
$ws.Range("C318:F318").Value = 6614488000000
$ws.Range("C319:F319").Value = 6708871000000
$ws.Range("C320:F320").Value = 6808405000000
$ws.Range("C327:F327").Value = 7965913000000
$ws.Range("C329:F329").Value = 8140532000000

# --- Append a new data row (330) ---
$ws.Range("A329").Copy()
$ws.Range("A330").PasteSpecial(-4122)
$ws.Range("A330").Value = 45078.41666666666
$ws.Range("B330").Value = "ECONOMICS:EGM2"
$ws.Range("C330").Value = 8248190000000
$ws.Range("D330").Value = 8248190000000
$ws.Range("E330").Value = 8248190000000
$ws.Range("F330").Value = 8248190000000
$ws.Range("G330").Value = 0
